# Commit: "changed ajax files with modified sql statements-Continued"
# Adds a new "Get list of previous prescriptions" query row (row 12) to the
# "Query" worksheet, together with its parameterised SQL (column C) and its
# worked example with literal values (column E), and moves the sheet's
# on-screen selection/scroll position down to show the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Query")

# --- New shared strings -----------------------------------------------

$label = "Get list of previous prescriptions"

$sqlParam = @'
SELECT a.VISIT_ID, a.PATIENT_ID, a.VISIT_DATE, b.prescription_id, a.visit_id 
                FROM visit a, prescription b
                WHERE a.patient_id = '$patient_id'
                AND a.visited = 'YES'
                AND a.visit_id = b.visit_id and a.chamber_id='".$chamber_name."' and a.doc_id='".$doc_name."' and a.chamber_id=b.chamber_id and a.doc_id=b.doc_id
                AND b.STATUS = 'SAVE' order by VISIT_DATE desc LIMIT 0 , 5
'@

$sqlExample = @'
SELECT a.VISIT_ID, a.PATIENT_ID, a.VISIT_DATE, b.prescription_id, a.visit_id 
                FROM visit a, prescription b
                WHERE a.patient_id = '123'
                AND a.visited = 'YES'
                AND a.visit_id = b.visit_id and a.chamber_id='sos' and a.doc_id='sroy' and a.chamber_id=b.chamber_id and a.doc_id=b.doc_id
                AND b.STATUS = 'SAVE' order by VISIT_DATE desc LIMIT 0 , 5
'@

# --- Row 12: Description / Complete Query / Example --------------------
# New shared strings are appended in the order the values are written, and
# the target workbook expects: 56 = label, 57 = parameterised SQL,
# 58 = worked example SQL - so populate D12 before C12/E12.

$ws.Range("D12").Value = $label
$ws.Range("C12").Value = $sqlParam
$ws.Range("E12").Value = $sqlExample

# Column C already carries the wrap-text style from the existing sheet;
# column E needs it applied explicitly (matches rows 9/10 above it).
$ws.Range("E12").WrapText = $true

# Row grows tall enough to show the multi-line query text.
$ws.Rows.Item(12).RowHeight = 105

# --- Scroll / selection housekeeping -----------------------------------

$ws.Activate()
$ws.Range("C13").Select()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
